$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(62, 8).Value = 97232824
$ws.Cells.Item(62, 9).Value = 45467800
$ws.Cells.Item(62, 10).Value = 178577860
$ws.Cells.Item(62, 11).Value = 45467800
$ws.Cells.Item(62, 12).Value = 178577860
$ws.Cells.Item(62, 13).Value = -45467176
$ws.Cells.Item(62, 14).Value = -178579108

$ws.Cells.Item(65, 8).Value = 97232824
$ws.Cells.Item(65, 9).Value = 45467800
$ws.Cells.Item(65, 10).Value = 178577860
$ws.Cells.Item(65, 11).Value = 227339000
$ws.Cells.Item(65, 12).Value = 892889300
$ws.Cells.Item(65, 13).Value = -227335880
$ws.Cells.Item(65, 14).Value = -892895540

$ws.Cells.Item(98, 8).Value = 71152456
$ws.Cells.Item(98, 9).Value = 33337000
$ws.Cells.Item(98, 10).Value = 127875640
$ws.Cells.Item(98, 11).Value = 33337000
$ws.Cells.Item(98, 12).Value = 127875640
$ws.Cells.Item(98, 13).Value = -33335502
$ws.Cells.Item(98, 14).Value = -127878636

$ws.Cells.Item(113, 8).Value = 8335361
$ws.Cells.Item(113, 9).Value = 10001800
$ws.Cells.Item(113, 10).Value = 3166.5
$ws.Cells.Item(113, 11).Value = 10001800
$ws.Cells.Item(113, 12).Value = 3166.5
$ws.Cells.Item(113, 13).Value = -9998546
$ws.Cells.Item(113, 14).Value = -9674.5

$ws.Cells.Item(122, 8).Value = 71152456
$ws.Cells.Item(122, 9).Value = 33337000
$ws.Cells.Item(122, 10).Value = 127875640
$ws.Cells.Item(122, 11).Value = 100011000
$ws.Cells.Item(122, 12).Value = 383626920
$ws.Cells.Item(122, 13).Value = -100008550
$ws.Cells.Item(122, 14).Value = -383631820

$ws.Cells.Item(129, 8).Value = 403343.47
$ws.Cells.Item(129, 9).Value = 2830.2307
$ws.Cells.Item(129, 10).Value = 837232.8
$ws.Cells.Item(129, 11).Value = 8490.6921
$ws.Cells.Item(129, 12).Value = 2511698.4
$ws.Cells.Item(129, 13).Value = -3490.6921
$ws.Cells.Item(129, 14).Value = -2521698.4

$ws.Cells.Item(132, 8).Value = 1793440.4
$ws.Cells.Item(132, 9).Value = 980.07275
$ws.Cells.Item(132, 11).Value = 2940.21825
$ws.Cells.Item(132, 13).Value = -410.2182499999999

$ws.Cells.Item(138, 8).Value = 2535.7605
$ws.Cells.Item(138, 9).Value = 1972.6285
$ws.Cells.Item(138, 10).Value = 3083.25
$ws.Cells.Item(138, 11).Value = 5917.8855
$ws.Cells.Item(138, 12).Value = 9249.75
$ws.Cells.Item(138, 13).Value = -777.8855000000003
$ws.Cells.Item(138, 14).Value = -19529.75

$ws.Cells.Item(141, 8).Value = 1216.1041
$ws.Cells.Item(141, 9).Value = 929.09753
$ws.Cells.Item(141, 10).Value = 2897.1428
$ws.Cells.Item(141, 11).Value = 2787.29259
$ws.Cells.Item(141, 12).Value = 8691.428400000001
$ws.Cells.Item(141, 13).Value = 2392.70741
$ws.Cells.Item(141, 14).Value = -19051.4284

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(61, 8).Value = 3350653.8
$ws.Cells.Item(61, 9).Value = 1489003
$ws.Cells.Item(61, 11).Value = 1489003
$ws.Cells.Item(61, 13).Value = -1488791

$ws.Cells.Item(132, 8).Value = 9418551
$ws.Cells.Item(132, 9).Value = 10641080
$ws.Cells.Item(132, 10).Value = 4630310.5
$ws.Cells.Item(132, 11).Value = 31923240
$ws.Cells.Item(132, 12).Value = 13890931.5
$ws.Cells.Item(132, 13).Value = -31920710
$ws.Cells.Item(132, 14).Value = -13895991.5

$ws.Cells.Item(136, 8).Value = 3350653.8
$ws.Cells.Item(136, 9).Value = 1489003
$ws.Cells.Item(136, 11).Value = 4467009
$ws.Cells.Item(136, 13).Value = -4464459

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(86, 8).Value = 1947.83
$ws.Cells.Item(86, 9).Value = 1961.5416
$ws.Cells.Item(86, 10).Value = 1618.75
$ws.Cells.Item(86, 11).Value = 1961.5416
$ws.Cells.Item(86, 12).Value = 1618.75
$ws.Cells.Item(86, 13).Value = -838.5416
$ws.Cells.Item(86, 14).Value = -3864.75

$ws.Cells.Item(89, 8).Value = 1947.83
$ws.Cells.Item(89, 9).Value = 1961.5416
$ws.Cells.Item(89, 10).Value = 1618.75
$ws.Cells.Item(89, 11).Value = 9807.708000000001
$ws.Cells.Item(89, 12).Value = 8093.75
$ws.Cells.Item(89, 13).Value = -4191.708000000001
$ws.Cells.Item(89, 14).Value = -19325.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 1456580.2
$ws.Cells.Item(31, 9).Value = 1154.2812
$ws.Cells.Item(31, 11).Value = 1154.2812
$ws.Cells.Item(31, 13).Value = -859.2811999999999

$ws.Cells.Item(34, 8).Value = 1456580.2
$ws.Cells.Item(34, 9).Value = 1154.2812
$ws.Cells.Item(34, 11).Value = 1154.2812
$ws.Cells.Item(34, 13).Value = -952.2811999999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(81, 8).Value = 6432.7085
$ws.Cells.Item(81, 9).Value = 1700
$ws.Cells.Item(81, 10).Value = 6862.9546
$ws.Cells.Item(81, 11).Value = 5100
$ws.Cells.Item(81, 12).Value = 20588.8638
$ws.Cells.Item(81, 13).Value = -3977
$ws.Cells.Item(81, 14).Value = -22834.8638

$ws.Cells.Item(84, 8).Value = 6432.7085
$ws.Cells.Item(84, 9).Value = 1700
$ws.Cells.Item(84, 10).Value = 6862.9546
$ws.Cells.Item(84, 11).Value = 15300
$ws.Cells.Item(84, 12).Value = 61766.5914
$ws.Cells.Item(84, 13).Value = -9684
$ws.Cells.Item(84, 14).Value = -72998.5914

$ws.Cells.Item(136, 8).Value = 1865.8334
$ws.Cells.Item(136, 9).Value = 1118.75
$ws.Cells.Item(136, 10).Value = 3360
$ws.Cells.Item(136, 11).Value = 3356.25
$ws.Cells.Item(136, 12).Value = 10080
$ws.Cells.Item(136, 13).Value = 1743.75
$ws.Cells.Item(136, 14).Value = -20280

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(102, 8).Value = 4013.4443
$ws.Cells.Item(102, 9).Value = 4308.875
$ws.Cells.Item(102, 10).Value = 1650
$ws.Cells.Item(102, 11).Value = 4308.875
$ws.Cells.Item(102, 12).Value = 1650
$ws.Cells.Item(102, 13).Value = -2686.875
$ws.Cells.Item(102, 14).Value = -4894

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 1564.7646
$ws.Cells.Item(7, 9).Value = 1483.2727
$ws.Cells.Item(7, 10).Value = 1714.1666
$ws.Cells.Item(7, 11).Value = 1483.2727
$ws.Cells.Item(7, 12).Value = 1714.1666
$ws.Cells.Item(7, 13).Value = -1371.2727
$ws.Cells.Item(7, 14).Value = -1938.1666

$ws.Cells.Item(40, 8).Value = 2817.4375
$ws.Cells.Item(40, 9).Value = 2954.4546
$ws.Cells.Item(40, 10).Value = 2516
$ws.Cells.Item(40, 11).Value = 2954.4546
$ws.Cells.Item(40, 12).Value = 2516
$ws.Cells.Item(40, 13).Value = -2818.4546
$ws.Cells.Item(40, 14).Value = -2788

$ws.Cells.Item(126, 8).Value = 1564.7646
$ws.Cells.Item(126, 9).Value = 1483.2727
$ws.Cells.Item(126, 10).Value = 1714.1666
$ws.Cells.Item(126, 11).Value = 4449.8181
$ws.Cells.Item(126, 12).Value = 5142.4998
$ws.Cells.Item(126, 13).Value = -1979.8181
$ws.Cells.Item(126, 14).Value = -10082.4998

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(107, 8).Value = 8034.231
$ws.Cells.Item(107, 9).Value = 10723.685
$ws.Cells.Item(107, 10).Value = 734.2857
$ws.Cells.Item(107, 11).Value = 32171.055
$ws.Cells.Item(107, 12).Value = 2202.8571
$ws.Cells.Item(107, 13).Value = -30251.055
$ws.Cells.Item(107, 14).Value = -6042.8571

$ws.Cells.Item(132, 8).Value = 1766056.2
$ws.Cells.Item(132, 9).Value = 2555.5454
$ws.Cells.Item(132, 10).Value = 9525460
$ws.Cells.Item(132, 11).Value = 7666.6362
$ws.Cells.Item(132, 12).Value = 28576380
$ws.Cells.Item(132, 13).Value = -5136.6362
$ws.Cells.Item(132, 14).Value = -28581440
